$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 420, pushing all existing rows (420..467) down to (421..468)
$ws.Rows("420:420").Insert()

# Populate the newly inserted row 420 with the new weekly data point
$ws.Range("A420").Value = 4
$ws.Range("B420").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C420").Value = "Los Lagos"
$ws.Range("D420").Value = 45194
$ws.Range("E420").Value = 10
$ws.Range("F420").Value = "Fruta"
$ws.Range("G420").Value = 100108
$ws.Range("H420").Value = "Tropicales y subtropicales"
$ws.Range("I420").Value = 100108005
$ws.Range("J420").Value = "Pi" + [char]241 + "a"
$ws.Range("K420").Value = "Caramelo"
$ws.Range("L420").Value = "Segunda"
$ws.Range("M420").Value = 100
$ws.Range("N420").Value = 25000
$ws.Range("O420").Value = 25000
$ws.Range("P420").Value = 25000
$ws.Range("Q420").Value = "$/caja 14 unidades"
$ws.Range("R420").Value = "Ecuador"
$ws.Range("S420").Value = 1786
$ws.Range("T420").Value = 14
